$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: set the "category" label for the header/first data row ---
$ws.Range("A2").Value = "Agriculture"

# --- Rows 81-201: fill in the category / sub_category sector names ---
$ws.Range("A81").Value = 'Energy_industry'
$ws.Range("B81").Value = 'Other_fuels'
$ws.Range("A82").Value = 'Energy_industry'
$ws.Range("B82").Value = 'Other_fuels'
$ws.Range("A83").Value = 'Energy_industry'
$ws.Range("B83").Value = 'Other_fuels'
$ws.Range("A84").Value = 'Energy_industry'
$ws.Range("B84").Value = 'Other_fuels'
$ws.Range("A85").Value = 'Energy_industry'
$ws.Range("B85").Value = 'Other_fuels'
$ws.Range("A86").Value = 'Energy_industry'
$ws.Range("B86").Value = 'Nuclear_fuel'
$ws.Range("A87").Value = 'Manufacturing_industry'
$ws.Range("B87").Value = 'Primary_materials'
$ws.Range("A88").Value = 'Manufacturing_industry'
$ws.Range("B88").Value = 'Secondary_materials'
$ws.Range("A89").Value = 'Manufacturing_industry'
$ws.Range("B89").Value = 'Other_manufacturing_industry'
$ws.Range("A90").Value = 'Manufacturing_industry'
$ws.Range("B90").Value = 'Other_manufacturing_industry'
$ws.Range("A91").Value = 'Manufacturing_industry'
$ws.Range("B91").Value = 'Other_manufacturing_industry'
$ws.Range("A92").Value = 'Energy_industry'
$ws.Range("B92").Value = 'Bioenergy'
$ws.Range("A93").Value = 'Energy_industry'
$ws.Range("B93").Value = 'Other_fuels'
$ws.Range("A94").Value = 'Energy_industry'
$ws.Range("B94").Value = 'Bioenergy'
$ws.Range("A95").Value = 'Energy_industry'
$ws.Range("B95").Value = 'Bioenergy'
$ws.Range("A96").Value = 'Energy_industry'
$ws.Range("B96").Value = 'Bioenergy'
$ws.Range("A97").Value = 'Manufacturing_industry'
$ws.Range("B97").Value = 'Materials'
$ws.Range("A98").Value = 'Manufacturing_industry'
$ws.Range("B98").Value = 'Primary_materials'
$ws.Range("A99").Value = 'Manufacturing_industry'
$ws.Range("B99").Value = 'Secondary_materials'
$ws.Range("A100").Value = 'Manufacturing_industry'
$ws.Range("B100").Value = 'Materials'
$ws.Range("A101").Value = 'Manufacturing_industry'
$ws.Range("B101").Value = 'Materials'
$ws.Range("A102").Value = 'Manufacturing_industry'
$ws.Range("B102").Value = 'Primary_materials'
$ws.Range("A103").Value = 'Manufacturing_industry'
$ws.Range("B103").Value = 'Secondary_materials'
$ws.Range("A104").Value = 'Manufacturing_industry'
$ws.Range("B104").Value = 'Materials'
$ws.Range("A105").Value = 'Manufacturing_industry'
$ws.Range("B105").Value = 'Primary_materials'
$ws.Range("A106").Value = 'Manufacturing_industry'
$ws.Range("B106").Value = 'Secondary_materials'
$ws.Range("A107").Value = 'Manufacturing_industry'
$ws.Range("B107").Value = 'Primary_materials'
$ws.Range("A108").Value = 'Manufacturing_industry'
$ws.Range("B108").Value = 'Secondary_materials'
$ws.Range("A109").Value = 'Manufacturing_industry'
$ws.Range("B109").Value = 'Primary_materials'
$ws.Range("A110").Value = 'Manufacturing_industry'
$ws.Range("B110").Value = 'Secondary_materials'
$ws.Range("A111").Value = 'Manufacturing_industry'
$ws.Range("B111").Value = 'Primary_materials'
$ws.Range("A112").Value = 'Manufacturing_industry'
$ws.Range("B112").Value = 'Secondary_materials'
$ws.Range("A113").Value = 'Manufacturing_industry'
$ws.Range("B113").Value = 'Primary_materials'
$ws.Range("A114").Value = 'Manufacturing_industry'
$ws.Range("B114").Value = 'Secondary_materials'
$ws.Range("A115").Value = 'Manufacturing_industry'
$ws.Range("B115").Value = 'Primary_materials'
$ws.Range("A116").Value = 'Manufacturing_industry'
$ws.Range("B116").Value = 'Secondary_materials'
$ws.Range("A117").Value = 'Manufacturing_industry'
$ws.Range("B117").Value = 'Other_manufacturing_industry'
$ws.Range("A118").Value = 'Manufacturing_industry'
$ws.Range("B118").Value = 'Other_manufacturing_industry'
$ws.Range("A119").Value = 'Manufacturing_industry'
$ws.Range("B119").Value = 'Other_manufacturing_industry'
$ws.Range("A120").Value = 'Manufacturing_industry'
$ws.Range("B120").Value = 'Other_manufacturing_industry'
$ws.Range("A121").Value = 'Manufacturing_industry'
$ws.Range("B121").Value = 'Other_manufacturing_industry'
$ws.Range("A122").Value = 'Manufacturing_industry'
$ws.Range("B122").Value = 'Other_manufacturing_industry'
$ws.Range("A123").Value = 'Manufacturing_industry'
$ws.Range("B123").Value = 'Other_manufacturing_industry'
$ws.Range("A124").Value = 'Manufacturing_industry'
$ws.Range("B124").Value = 'Transport_equipment'
$ws.Range("A125").Value = 'Manufacturing_industry'
$ws.Range("B125").Value = 'Transport_equipment'
$ws.Range("A126").Value = 'Manufacturing_industry'
$ws.Range("B126").Value = 'Other_manufacturing_industry'
$ws.Range("A127").Value = 'Manufacturing_industry'
$ws.Range("B127").Value = 'Other_manufacturing_industry'
$ws.Range("A128").Value = 'Manufacturing_industry'
$ws.Range("B128").Value = 'Other_manufacturing_industry'
$ws.Range("A129").Value = 'Energy_industry'
$ws.Range("B129").Value = 'Electricity'
$ws.Range("A130").Value = 'Energy_industry'
$ws.Range("B130").Value = 'Electricity'
$ws.Range("A131").Value = 'Energy_industry'
$ws.Range("B131").Value = 'Electricity'
$ws.Range("A132").Value = 'Energy_industry'
$ws.Range("B132").Value = 'Electricity'
$ws.Range("A133").Value = 'Energy_industry'
$ws.Range("B133").Value = 'Electricity'
$ws.Range("A134").Value = 'Energy_industry'
$ws.Range("B134").Value = 'Electricity'
$ws.Range("A135").Value = 'Energy_industry'
$ws.Range("B135").Value = 'Electricity'
$ws.Range("A136").Value = 'Energy_industry'
$ws.Range("B136").Value = 'Electricity'
$ws.Range("A137").Value = 'Energy_industry'
$ws.Range("B137").Value = 'Electricity'
$ws.Range("A138").Value = 'Energy_industry'
$ws.Range("B138").Value = 'Electricity'
$ws.Range("A139").Value = 'Energy_industry'
$ws.Range("B139").Value = 'Electricity'
$ws.Range("A140").Value = 'Energy_industry'
$ws.Range("B140").Value = 'Electricity'
$ws.Range("A141").Value = 'Energy_industry'
$ws.Range("B141").Value = 'Energy_Services'
$ws.Range("A142").Value = 'Energy_industry'
$ws.Range("B142").Value = 'Energy_Services'
$ws.Range("A143").Value = 'Energy_industry'
$ws.Range("B143").Value = 'Gas'
$ws.Range("A144").Value = 'Energy_industry'
$ws.Range("B144").Value = 'Gas'
$ws.Range("A145").Value = 'Energy_industry'
$ws.Range("B145").Value = 'Gas'
$ws.Range("A146").Value = 'Energy_industry'
$ws.Range("B146").Value = 'Gas'
$ws.Range("A147").Value = 'Energy_industry'
$ws.Range("B147").Value = 'Bioenergy'
$ws.Range("A148").Value = 'Energy_industry'
$ws.Range("B148").Value = 'Energy_Services'
$ws.Range("A149").Value = 'Energy_industry'
$ws.Range("B149").Value = 'Heat'
$ws.Range("A150").Value = 'Water_and_waste_treatment'
$ws.Range("B150").Value = 'Water_services'
$ws.Range("A151").Value = 'Manufacturing_industry'
$ws.Range("B151").Value = 'Primary_materials'
$ws.Range("A152").Value = 'Manufacturing_industry'
$ws.Range("B152").Value = 'Secondary_materials'
$ws.Range("A153").Value = 'Oth_services'
$ws.Range("B153").Value = 'Business_services'
$ws.Range("A154").Value = 'Oth_services'
$ws.Range("B154").Value = 'Trade_services'
$ws.Range("A155").Value = 'Oth_services'
$ws.Range("B155").Value = 'Trade_services'
$ws.Range("A156").Value = 'Oth_services'
$ws.Range("B156").Value = 'Trade_services'
$ws.Range("A157").Value = 'Oth_services'
$ws.Range("B157").Value = 'Business_services'
$ws.Range("A158").Value = 'Transports'
$ws.Range("B158").Value = 'Land_transport'
$ws.Range("A159").Value = 'Transports'
$ws.Range("B159").Value = 'Land_transport'
$ws.Range("A160").Value = 'Transports'
$ws.Range("B160").Value = 'Land_transport'
$ws.Range("A161").Value = 'Transports'
$ws.Range("B161").Value = 'Water_transport'
$ws.Range("A162").Value = 'Transports'
$ws.Range("B162").Value = 'Water_transport'
$ws.Range("A163").Value = 'Transports'
$ws.Range("B163").Value = 'Air_transport'
$ws.Range("A164").Value = 'Oth_services'
$ws.Range("B164").Value = 'Business_services'
$ws.Range("A165").Value = 'Oth_services'
$ws.Range("B165").Value = 'Business_services'
$ws.Range("A166").Value = 'Oth_services'
$ws.Range("B166").Value = 'Business_services'
$ws.Range("A167").Value = 'Oth_services'
$ws.Range("B167").Value = 'Business_services'
$ws.Range("A168").Value = 'Oth_services'
$ws.Range("B168").Value = 'Business_services'
$ws.Range("A169").Value = 'Oth_services'
$ws.Range("B169").Value = 'Business_services'
$ws.Range("A170").Value = 'Oth_services'
$ws.Range("B170").Value = 'Business_services'
$ws.Range("A171").Value = 'Oth_services'
$ws.Range("B171").Value = 'Business_services'
$ws.Range("A172").Value = 'Oth_services'
$ws.Range("B172").Value = 'Business_services'
$ws.Range("A173").Value = 'Oth_services'
$ws.Range("B173").Value = 'Business_services'
$ws.Range("A174").Value = 'Oth_services'
$ws.Range("B174").Value = 'Public_services'
$ws.Range("A175").Value = 'Oth_services'
$ws.Range("B175").Value = 'Public_services'
$ws.Range("A176").Value = 'Oth_services'
$ws.Range("B176").Value = 'Public_services'
$ws.Range("A177").Value = 'Water_and_waste_treatment'
$ws.Range("B177").Value = 'Incineration'
$ws.Range("A178").Value = 'Water_and_waste_treatment'
$ws.Range("B178").Value = 'Incineration'
$ws.Range("A179").Value = 'Water_and_waste_treatment'
$ws.Range("B179").Value = 'Incineration'
$ws.Range("A180").Value = 'Water_and_waste_treatment'
$ws.Range("B180").Value = 'Incineration'
$ws.Range("A181").Value = 'Water_and_waste_treatment'
$ws.Range("B181").Value = 'Incineration'
$ws.Range("A182").Value = 'Water_and_waste_treatment'
$ws.Range("B182").Value = 'Incineration'
$ws.Range("A183").Value = 'Water_and_waste_treatment'
$ws.Range("B183").Value = 'Incineration'
$ws.Range("A184").Value = 'Water_and_waste_treatment'
$ws.Range("B184").Value = 'Land_application'
$ws.Range("A185").Value = 'Water_and_waste_treatment'
$ws.Range("B185").Value = 'Land_application'
$ws.Range("A186").Value = 'Water_and_waste_treatment'
$ws.Range("B186").Value = 'Land_application'
$ws.Range("A187").Value = 'Water_and_waste_treatment'
$ws.Range("B187").Value = 'Land_application'
$ws.Range("A188").Value = 'Water_and_waste_treatment'
$ws.Range("B188").Value = 'Land_application'
$ws.Range("A189").Value = 'Water_and_waste_treatment'
$ws.Range("B189").Value = 'Other_waste_treatment'
$ws.Range("A190").Value = 'Water_and_waste_treatment'
$ws.Range("B190").Value = 'Other_waste_treatment'
$ws.Range("A191").Value = 'Water_and_waste_treatment'
$ws.Range("B191").Value = 'Landfill'
$ws.Range("A192").Value = 'Water_and_waste_treatment'
$ws.Range("B192").Value = 'Landfill'
$ws.Range("A193").Value = 'Water_and_waste_treatment'
$ws.Range("B193").Value = 'Landfill'
$ws.Range("A194").Value = 'Water_and_waste_treatment'
$ws.Range("B194").Value = 'Landfill'
$ws.Range("A195").Value = 'Water_and_waste_treatment'
$ws.Range("B195").Value = 'Landfill'
$ws.Range("A196").Value = 'Water_and_waste_treatment'
$ws.Range("B196").Value = 'Landfill'
$ws.Range("A197").Value = 'Oth_services'
$ws.Range("B197").Value = 'Business_services'
$ws.Range("A198").Value = 'Oth_services'
$ws.Range("B198").Value = 'Business_services'
$ws.Range("A199").Value = 'Oth_services'
$ws.Range("B199").Value = 'Business_services'
$ws.Range("A200").Value = 'Oth_services'
$ws.Range("B200").Value = 'Business_services'
$ws.Range("A201").Value = 'Oth_services'
$ws.Range("B201").Value = 'Business_services'

# --- Column layout: split former A:B (same width) so column A is wider
#     to fit the new longer category names, leaving column B untouched ---
$ws.Columns.Item(1).ColumnWidth = 14.83

# --- sheet view: scroll/select around the area that was being edited ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 103
$ws.Range("B117").Select()

# --- Conditional-format scratch styles (author previewed/applied a
#     cyan highlight + Arial font while aggregating the sectors, then
#     removed the rule, leaving the differential styles behind) ---
for ($i = 0; $i -lt 6; $i++) {
    $tmpRange = $ws.Range("A1")
    $fc = $tmpRange.FormatConditions.Add(1, 3, "1")
    $fc.Font.Name = "Arial"
    $fc.Interior.Color = 16777062
    $fc.Delete()
}
